# Handback status report regeneration:
#  - "Latest HO Xliff Generate Date" (Overview sheet) refreshed
#  - Priority flipped from human translation ("ht") to machine translation ("mt")
#  - Handoff / Handback datetimes refreshed to the new report run

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 20:20:42"
$wsOverview.Range("G3").Value = "2016-09-05 20:20:42"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-05 20:20:36"
$wsZhCn.Range("H3").Value = "2016-09-05 20:20:36"
$wsZhCn.Range("K2").Value = "2016-09-05 20:20:54"
$wsZhCn.Range("K3").Value = "2016-09-05 20:20:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-05 20:20:42"
$wsDeDe.Range("H3").Value = "2016-09-05 20:20:42"
$wsDeDe.Range("K2").Value = "2016-09-05 20:21:08"
$wsDeDe.Range("K3").Value = "2016-09-05 20:21:08"
